$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("33").Delete()
